# Update the "Förändrad" (Changed) date column (C) for rows 2-10 from
# 2023-10-22 (45221) to 2023-10-25 (45224).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = (Get-Date -Year 2023 -Month 10 -Day 25).Date

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
